$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 608.2381
$ws.Range("I28").Value = 613.65
$ws.Range("J28").Value = 500
$ws.Range("K28").Value = 613.65
$ws.Range("L28").Value = 500
$ws.Range("M28").Value = -128.65
$ws.Range("N28").Value = -1470
$ws.Range("H41").Value = 628.3570999999999
$ws.Range("I41").Value = 98.5
$ws.Range("J41").Value = 716.6667
$ws.Range("K41").Value = 98.5
$ws.Range("L41").Value = 716.6667
$ws.Range("M41").Value = 341.5
$ws.Range("N41").Value = -1596.6667
$ws.Range("H62").Value = 4125
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 5400
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 5400
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -6648
$ws.Range("H65").Value = 4125
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 5400
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 27000
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -33240
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = ""
$ws.Range("H76").Value = 4000
$ws.Range("I76").Value = 4000
$ws.Range("K76").Value = 4000
$ws.Range("M76").Value = -3685
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = ""
$ws.Range("H79").Value = 4000
$ws.Range("I79").Value = 4000
$ws.Range("K79").Value = 4000
$ws.Range("M79").Value = -2908
$ws.Range("H86").Value = 7583.8335
$ws.Range("I86").Value = 1751.5
$ws.Range("K86").Value = 1751.5
$ws.Range("M86").Value = -628.5
$ws.Range("H89").Value = 7583.8335
$ws.Range("I89").Value = 1751.5
$ws.Range("K89").Value = 8757.5
$ws.Range("M89").Value = -3141.5
$ws.Range("H106").Value = 4166.6665
$ws.Range("I106").Value = 2500
$ws.Range("K106").Value = 2500
$ws.Range("M106").Value = -1869
$ws.Range("H116").Value = 14430.6
$ws.Range("I116").Value = 4833.3335
$ws.Range("J116").Value = 18543.715
$ws.Range("K116").Value = 4833.3335
$ws.Range("L116").Value = 18543.715
$ws.Range("M116").Value = -1391.3335
$ws.Range("N116").Value = -25427.715
$ws.Range("H129").Value = 982.8778
$ws.Range("J129").Value = 1024.2235
$ws.Range("L129").Value = 3072.6705
$ws.Range("N129").Value = -13072.6705

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1643.6666
$ws.Range("I2").Value = 1643.6666
$ws.Range("K2").Value = 1643.6666
$ws.Range("M2").Value = -1530.6666
$ws.Range("H61").Value = 2088
$ws.Range("I61").Value = 1672.2667
$ws.Range("K61").Value = 1672.2667
$ws.Range("M61").Value = -1460.2667
$ws.Range("H109").Value = 28000
$ws.Range("J109").Value = 28000
$ws.Range("L109").Value = 28000
$ws.Range("N109").Value = -30774
$ws.Range("H116").Value = 1643.6666
$ws.Range("I116").Value = 1643.6666
$ws.Range("K116").Value = 1643.6666
$ws.Range("M116").Value = 650.3334
$ws.Range("H122").Value = 2094.6072
$ws.Range("I122").Value = 1260.5
$ws.Range("J122").Value = 3596
$ws.Range("K122").Value = 3781.5
$ws.Range("L122").Value = 10788
$ws.Range("M122").Value = -1331.5
$ws.Range("N122").Value = -15688
$ws.Range("H136").Value = 2088
$ws.Range("I136").Value = 1672.2667
$ws.Range("K136").Value = 5016.800099999999
$ws.Range("M136").Value = -2466.800099999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1643.6666
$ws.Range("I3").Value = 1643.6666
$ws.Range("K3").Value = 1643.6666
$ws.Range("M3").Value = -1529.6666
$ws.Range("H75").Value = 12526.889
$ws.Range("I75").Value = 9435.5
$ws.Range("J75").Value = 15000
$ws.Range("K75").Value = 9435.5
$ws.Range("L75").Value = 15000
$ws.Range("M75").Value = -8499.5
$ws.Range("N75").Value = -16872
$ws.Range("H78").Value = 12526.889
$ws.Range("I78").Value = 9435.5
$ws.Range("J78").Value = 15000
$ws.Range("K78").Value = 28306.5
$ws.Range("L78").Value = 45000
$ws.Range("M78").Value = -23626.5
$ws.Range("N78").Value = -54360
$ws.Range("H92").Value = 60000
$ws.Range("J92").Value = 60000
$ws.Range("L92").Value = 60000
$ws.Range("N92").Value = -64992
$ws.Range("H94").Value = 1936.2354
$ws.Range("I94").Value = 1453.1538
$ws.Range("J94").Value = 3506.25
$ws.Range("K94").Value = 1453.1538
$ws.Range("L94").Value = 3506.25
$ws.Range("M94").Value = -1002.1538
$ws.Range("N94").Value = -4408.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 15879021
$ws.Range("I16").Value = 37041708
$ws.Range("J16").Value = 7006
$ws.Range("K16").Value = 37041708
$ws.Range("L16").Value = 7006
$ws.Range("M16").Value = -37041421
$ws.Range("N16").Value = -7580
$ws.Range("H31").Value = 11113817
$ws.Range("I31").Value = 1392.7097
$ws.Range("K31").Value = 1392.7097
$ws.Range("M31").Value = -1097.7097
$ws.Range("H34").Value = 11113817
$ws.Range("I34").Value = 1392.7097
$ws.Range("K34").Value = 1392.7097
$ws.Range("M34").Value = -1190.7097
$ws.Range("H58").Value = 1904.8616
$ws.Range("I58").Value = 1649.8036
$ws.Range("J58").Value = 3491.889
$ws.Range("K58").Value = 1649.8036
$ws.Range("L58").Value = 3491.889
$ws.Range("M58").Value = -1446.8036
$ws.Range("N58").Value = -3897.889
$ws.Range("H94").Value = 1328.6666
$ws.Range("J94").Value = 1478
$ws.Range("L94").Value = 1478
$ws.Range("N94").Value = -2380
$ws.Range("H107").Value = 902.53845
$ws.Range("I107").Value = 647.8570999999999
$ws.Range("K107").Value = 647.8570999999999
$ws.Range("M107").Value = 1272.1429
$ws.Range("H113").Value = 15879021
$ws.Range("I113").Value = 37041708
$ws.Range("J113").Value = 7006
$ws.Range("K113").Value = 37041708
$ws.Range("L113").Value = 7006
$ws.Range("M113").Value = -37039538
$ws.Range("N113").Value = -11346
$ws.Range("H134").Value = 4027.5813
$ws.Range("I134").Value = 4316.6206
$ws.Range("K134").Value = 12949.8618
$ws.Range("M134").Value = -10414.8618
$ws.Range("H136").Value = 1904.8616
$ws.Range("I136").Value = 1649.8036
$ws.Range("J136").Value = 3491.889
$ws.Range("K136").Value = 4949.4108
$ws.Range("L136").Value = 10475.667
$ws.Range("M136").Value = -2399.4108
$ws.Range("N136").Value = -15575.667

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 3750
$ws.Range("I87").Value = 3750
$ws.Range("K87").Value = 11250
$ws.Range("M87").Value = -10002
$ws.Range("H90").Value = 3750
$ws.Range("I90").Value = 3750
$ws.Range("K90").Value = 33750
$ws.Range("M90").Value = -27510

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4165.9414
$ws.Range("I122").Value = 2364.25
$ws.Range("K122").Value = 7092.75
$ws.Range("M122").Value = -4642.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1782
$ws.Range("I16").Value = 1536.4286
$ws.Range("J16").Value = 3501
$ws.Range("K16").Value = 1536.4286
$ws.Range("L16").Value = 3501
$ws.Range("M16").Value = -1366.4286
$ws.Range("N16").Value = -3841
$ws.Range("H22").Value = 2730.9375
$ws.Range("I22").Value = 1950.25
$ws.Range("J22").Value = 2991.1667
$ws.Range("K22").Value = 1950.25
$ws.Range("L22").Value = 2991.1667
$ws.Range("M22").Value = -1655.25
$ws.Range("N22").Value = -3581.1667
$ws.Range("H27").Value = 2730.9375
$ws.Range("I27").Value = 1950.25
$ws.Range("J27").Value = 2991.1667
$ws.Range("K27").Value = 1950.25
$ws.Range("L27").Value = 2991.1667
$ws.Range("M27").Value = -1843.25
$ws.Range("N27").Value = -3205.1667
$ws.Range("H122").Value = 4966.7896
$ws.Range("I122").Value = 2848.4
$ws.Range("J122").Value = 7320.5557
$ws.Range("K122").Value = 8545.200000000001
$ws.Range("L122").Value = 21961.6671
$ws.Range("M122").Value = -6095.200000000001
$ws.Range("N122").Value = -26861.6671
$ws.Range("H132").Value = 2987.8333
$ws.Range("I132").Value = 1548.2667
$ws.Range("J132").Value = 7306.533
$ws.Range("K132").Value = 4644.800099999999
$ws.Range("L132").Value = 21919.599
$ws.Range("M132").Value = -2114.800099999999
$ws.Range("N132").Value = -26979.599

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 40344
$ws.Range("J95").Value = 40344
$ws.Range("L95").Value = 40344
$ws.Range("N95").Value = -45836
$ws.Range("H113").Value = 2963.2632
$ws.Range("I113").Value = 4094.8845
$ws.Range("K113").Value = 12284.6535
$ws.Range("M113").Value = -10114.6535
$ws.Range("H122").Value = 3055.4666
$ws.Range("I122").Value = 1638.2
$ws.Range("K122").Value = 4914.6
$ws.Range("M122").Value = -2464.6
$ws.Range("H123").Value = 35520
$ws.Range("J123").Value = 35520
$ws.Range("L123").Value = 35520
$ws.Range("N123").Value = -45320
$ws.Range("H132").Value = 12826340
$ws.Range("I132").Value = 7633.7334
$ws.Range("J132").Value = 30306396
$ws.Range("K132").Value = 22901.2002
$ws.Range("L132").Value = 90919188
$ws.Range("M132").Value = -20371.2002
$ws.Range("N132").Value = -90924248
$ws.Range("H136").Value = 1451.0625
$ws.Range("I136").Value = 833.8570999999999
$ws.Range("J136").Value = 2629.3635
$ws.Range("K136").Value = 2501.5713
$ws.Range("L136").Value = 7888.0905
$ws.Range("M136").Value = 48.42870000000039
$ws.Range("N136").Value = -12988.0905
